$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (Dr. Wilford Rempel): last_name column (B) should only contain the last name "Rempel"
$ws.Range("B9").Value = "Rempel"

# Row 10 (Ms. Damaris Luettgen MD): last_name column (B) should contain "Damaris Luettgen MD"
$ws.Range("B10").Value = "Damaris Luettgen MD"

# Update the active cell/selection to B12 (next empty row under the imported data)
$ws.Range("B12").Select()
